$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Compra" -> "Compras"
# ---------------------------------------------------------------------------
$wsCompras = $wb.Worksheets.Item("Compra")
$wsCompras.Name = "Compras"

# ---------------------------------------------------------------------------
# 2. Refresh the remembered selection on the first four sheets (cosmetic,
#    matches the saved cursor position captured the next time the workbook
#    was opened/saved).
# ---------------------------------------------------------------------------
$wsCompras.Range("A2").Select()

$wsSector = $wb.Worksheets.Item("Sector_venta")
$wsSector.Range("A3").Select()

$wsVendedor = $wb.Worksheets.Item("Vendedor")
$wsVendedor.Range("A3").Select()

$wsFecha = $wb.Worksheets.Item("Fecha")
$wsFecha.Range("A4").Select()

$wsProducto = $wb.Worksheets.Item("Producto")
$wsProducto.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Rebuild the "Venta" sheet: five new leading columns (id_lugar_compra,
#    id_sector_venta, id_vendedor, id_fecha, id_producto) are added in front
#    of the existing num_nota_venta..Kilos table, which now lives in F:M
#    instead of A:H. Cell values are written directly into their final
#    positions (rather than using an Insert, which would also drag the
#    col-A width formatting over to col F) so the <cols> width entry stays
#    anchored on column A exactly as in the saved workbook.
# ---------------------------------------------------------------------------
$wsVenta = $wb.Worksheets.Item("Venta")

# Headers
$wsVenta.Cells.Item(1,1).Value = "id_lugar_compra"
$wsVenta.Cells.Item(1,2).Value = "id_sector_venta"
$wsVenta.Cells.Item(1,3).Value = "id_vendedor"
$wsVenta.Cells.Item(1,4).Value = "id_fecha"
$wsVenta.Cells.Item(1,5).Value = "id_producto"
$wsVenta.Cells.Item(1,6).Value = "num_nota_venta"
$wsVenta.Cells.Item(1,7).Value = "not_cod_cd"
$wsVenta.Cells.Item(1,8).Value = "MontoVenta"
$wsVenta.Cells.Item(1,9).Value = "MontoCosto"
$wsVenta.Cells.Item(1,10).Value = "Ganancia"
$wsVenta.Cells.Item(1,11).Value = "Unidades"
$wsVenta.Cells.Item(1,12).Value = "Volumen"
$wsVenta.Cells.Item(1,13).Value = "Kilos"

# Row 2
$wsVenta.Cells.Item(2,1).Value = 1
$wsVenta.Cells.Item(2,2).Value = 2
$wsVenta.Cells.Item(2,3).Value = 2
$wsVenta.Cells.Item(2,4).Value = 3
$wsVenta.Cells.Item(2,5).Value = 9
$wsVenta.Cells.Item(2,6).Value = 1
$wsVenta.Cells.Item(2,7).Value = 1
$wsVenta.Cells.Item(2,8).Value = 2000
$wsVenta.Cells.Item(2,9).Value = 1350
$wsVenta.Cells.Item(2,10).Value = 650
$wsVenta.Cells.Item(2,11).Value = 300
$wsVenta.Cells.Item(2,12).Value = 300
$wsVenta.Cells.Item(2,13).Value = 300

# Row 3
$wsVenta.Cells.Item(3,1).Value = 2
$wsVenta.Cells.Item(3,2).Value = 1
$wsVenta.Cells.Item(3,3).Value = 4
$wsVenta.Cells.Item(3,4).Value = 6
$wsVenta.Cells.Item(3,5).Value = 10
$wsVenta.Cells.Item(3,6).Value = 2
$wsVenta.Cells.Item(3,7).Value = 2
$wsVenta.Cells.Item(3,8).Value = 1500
$wsVenta.Cells.Item(3,9).Value = 1000
$wsVenta.Cells.Item(3,10).Value = 500
$wsVenta.Cells.Item(3,11).Value = 200
$wsVenta.Cells.Item(3,12).Value = 200
$wsVenta.Cells.Item(3,13).Value = 200

# Row 4
$wsVenta.Cells.Item(4,1).Value = 3
$wsVenta.Cells.Item(4,2).Value = 3
$wsVenta.Cells.Item(4,3).Value = 1
$wsVenta.Cells.Item(4,4).Value = 2
$wsVenta.Cells.Item(4,5).Value = 7
$wsVenta.Cells.Item(4,6).Value = 3
$wsVenta.Cells.Item(4,7).Value = 3
$wsVenta.Cells.Item(4,8).Value = 2350
$wsVenta.Cells.Item(4,9).Value = 1890
$wsVenta.Cells.Item(4,10).Value = 460
$wsVenta.Cells.Item(4,11).Value = 460
$wsVenta.Cells.Item(4,12).Value = 460
$wsVenta.Cells.Item(4,13).Value = 260

# Row 5
$wsVenta.Cells.Item(5,1).Value = 4
$wsVenta.Cells.Item(5,2).Value = 6
$wsVenta.Cells.Item(5,3).Value = 5
$wsVenta.Cells.Item(5,4).Value = 1
$wsVenta.Cells.Item(5,5).Value = 6
$wsVenta.Cells.Item(5,6).Value = 4
$wsVenta.Cells.Item(5,7).Value = 4
$wsVenta.Cells.Item(5,8).Value = 785
$wsVenta.Cells.Item(5,9).Value = 490
$wsVenta.Cells.Item(5,10).Value = 295
$wsVenta.Cells.Item(5,11).Value = 135
$wsVenta.Cells.Item(5,12).Value = 135
$wsVenta.Cells.Item(5,13).Value = 135

# Row 6
$wsVenta.Cells.Item(6,1).Value = 5
$wsVenta.Cells.Item(6,2).Value = 8
$wsVenta.Cells.Item(6,3).Value = 7
$wsVenta.Cells.Item(6,4).Value = 4
$wsVenta.Cells.Item(6,5).Value = 4
$wsVenta.Cells.Item(6,6).Value = 5
$wsVenta.Cells.Item(6,7).Value = 5
$wsVenta.Cells.Item(6,8).Value = 978
$wsVenta.Cells.Item(6,9).Value = 580
$wsVenta.Cells.Item(6,10).Value = 398
$wsVenta.Cells.Item(6,11).Value = 160
$wsVenta.Cells.Item(6,12).Value = 160
$wsVenta.Cells.Item(6,13).Value = 160

# Column widths for the two new columns that picked up a custom width;
# column A's original custom width (15.44140625) is left untouched.
$wsVenta.Columns("B").ColumnWidth = 17.21875
$wsVenta.Columns("C").ColumnWidth = 11.109375

# Final cursor position / active sheet, matching the saved workbook state.
$wsVenta.Range("F15").Select()
